$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column C (Fitness) values for rows 2-166 per new log/run data
$ranges = @(
    @{ StartRow = 2; EndRow = 2; Value = 12979 }
    @{ StartRow = 3; EndRow = 8; Value = 10261 }
    @{ StartRow = 9; EndRow = 12; Value = 9745 }
    @{ StartRow = 13; EndRow = 19; Value = 8666 }
    @{ StartRow = 20; EndRow = 121; Value = 7884 }
    @{ StartRow = 122; EndRow = 123; Value = 7345 }
    @{ StartRow = 124; EndRow = 126; Value = 7320 }
    @{ StartRow = 127; EndRow = 154; Value = 7312 }
    @{ StartRow = 155; EndRow = 166; Value = 7295 }
)

foreach ($r in $ranges) {
    for ($row = $r.StartRow; $row -le $r.EndRow; $row++) {
        $ws.Cells.Item($row, 3).Value2 = $r.Value
    }
}

Write-Output "Done updating C2:C166"
